# Insert a new weekly price record as the first (most recent) row of the
# "Femacal de La Calera - Pepino ensalada" block, pushing the existing
# rows 379..472 down to 380..473.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(379).Insert()

$ws.Cells.Item(379, 1).Value  = 3
$ws.Cells.Item(379, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(379, 3).Value  = "Coquimbo"
$ws.Cells.Item(379, 4).Value  = 44932
$ws.Cells.Item(379, 5).Value  = 5
$ws.Cells.Item(379, 6).Value  = 100112043
$ws.Cells.Item(379, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(379, 8).Value  = "Sin especificar"
$ws.Cells.Item(379, 9).Value  = "Primera"
$ws.Cells.Item(379, 10).Value = 85
$ws.Cells.Item(379, 11).Value = 16000
$ws.Cells.Item(379, 12).Value = 17000
$ws.Cells.Item(379, 13).Value = 16471
$ws.Cells.Item(379, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(379, 15).Value = "Limache"
$ws.Cells.Item(379, 16).Value = 275
$ws.Cells.Item(379, 17).Value = 60
$ws.Cells.Item(379, 18).Value = "Hortaliza"
